# Minor reformatting of several slides (Records deck)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 16 - "Constructor for Class RecordType"
# ---------------------------------------------------------------------------
$s16 = $p.Slides.Item(16)

# Shape 5: "TextBox 5" - widen/move the note box and fix its wording.
$noteBox = $s16.Shapes.Item(5)
$noteBox.Left  = 216
$noteBox.Width = 399.1351318359375

$noteText = $noteBox.TextFrame.TextRange
$fullNote = $noteText.Text
$idx = $fullNote.IndexOf("lambda expression to compute size.")
$oldLambda = "lambda expression to compute size."
$sub = $noteText.Characters($idx + 1, $oldLambda.Length)
$sub.Text = "lambda expression to compute record size."

# Shape 6: "Diamond 6" - nudge right slightly.
$diamond = $s16.Shapes.Item(6)
$diamond.Left = 408.3675842285156

# Shape 7: "Connector: Elbow 8" - re-route alongside the diamond; it is no
# longer horizontally flipped and collapses to a zero-width vertical line.
$connector = $s16.Shapes.Item(7)
$connector.HorizontalFlip = 0
$connector.Left  = 415.56756591796875
$connector.Width = 0

# ---------------------------------------------------------------------------
# Slide 19 - shift the textbox over the code listing
# ---------------------------------------------------------------------------
$s19 = $p.Slides.Item(19)
$tb19 = $s19.Shapes.Item(5)
$tb19.Left = 109.25

# ---------------------------------------------------------------------------
# Slide 21 - "record" example: tighten a comment's leading whitespace and
# add a little breathing room above it.
# ---------------------------------------------------------------------------
$s21 = $p.Slides.Item(21)
$codeBox = $s21.Shapes.Item(2)
$codeRange = $codeBox.TextFrame.TextRange
$fullCode = $codeRange.Text

$oldRectLine = "type Rectangle = record                 // fields are records"
$newRectLine = "type Rectangle = record        // fields are records"
$rectIdx = $fullCode.IndexOf($oldRectLine)
$rectRange = $codeRange.Characters($rectIdx + 1, $oldRectLine.Length)
$rectRange.Text = $newRectLine

# The blank paragraph immediately preceding "type Rectangle ..." gets a
# little extra space before it. Target it via the two paragraph-mark
# characters that straddle that (otherwise empty) paragraph.
$blankParaPos = $rectIdx - 1
$blankPara = $codeRange.Characters($blankParaPos, 2)
$blankPara.ParagraphFormat.SpaceBefore = 5

# ---------------------------------------------------------------------------
# Slide 25 - raise the textbox a touch
# ---------------------------------------------------------------------------
$s25 = $p.Slides.Item(25)
$tb25 = $s25.Shapes.Item(5)
$tb25.Top = 417.02764892578125
